$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2013 exam rows
$ws.Range("A11").Value = "2013 - Vår"
$ws.Range("B11").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-13-v.pdf)"
$ws.Range("C11").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-13-v-fasit.pdf)"

$ws.Range("A12").Value = "2013 - Høst"
$ws.Range("B12").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-13-h.pdf)"
$ws.Range("C12").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-13-h-fasit.pdf)"
$ws.Range("D12").Value = "1g-1h"

# Shorten the old task-reference notes in column D (drop the "Oppgave " prefix),
# working back up from the most recent row
$ws.Range("D10").Value = "1a (ANOVA-delen)"
$ws.Range("D8").Value = "1a-1f, 3c-3e"
$ws.Range("D7").Value = "1d-1e"
$ws.Range("D6").Value = "3a-3c"
$ws.Range("D5").Value = "1c, 3a-3e"
$ws.Range("D3").Value = "1, 2c-2e"
$ws.Range("D2").Value = "1c-1e, 3a-3c, 4a-4c"

# Add the remaining exam rows (2014 - 2020)
$ws.Range("A13").Value = "2014 - Vår"
$ws.Range("B13").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-14-v.pdf)"
$ws.Range("C13").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-14-v-fasit.pdf)"
$ws.Range("D13").Value = "1f-1i"

$ws.Range("A14").Value = "2014 - Høst"
$ws.Range("B14").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-14-h.pdf)"
$ws.Range("C14").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-14-h-fasit.pdf)"

$ws.Range("A15").Value = "2015 - Vår"
$ws.Range("B15").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-15-v.pdf)"
$ws.Range("C15").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-15-v-fasit.pdf)"

$ws.Range("A16").Value = "2015 - Høst"
$ws.Range("B16").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-15-h.pdf)"
$ws.Range("C16").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-15-h-fasit.pdf)"
$ws.Range("D16").Value = "3a-3c"

$ws.Range("A17").Value = "2016 - Vår"
$ws.Range("B17").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-16-v.pdf)"
$ws.Range("C17").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-16-v-fasit.pdf)"

$ws.Range("A18").Value = "2016 - Høst"
$ws.Range("B18").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-16-h.pdf)"
$ws.Range("C18").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-16-h-fasit.pdf)"

$ws.Range("A19").Value = "2017 - Vår"
$ws.Range("B19").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-17-v.pdf)"
$ws.Range("C19").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-17-v-fasit.pdf)"
$ws.Range("D19").Value = "1c"

$ws.Range("A20").Value = "2017 - Høst"
$ws.Range("B20").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-17-h.pdf)"
$ws.Range("C20").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-17-h-fasit.pdf)"
$ws.Range("D20").Value = "1f-1g "

$ws.Range("A21").Value = "2018 - Vår"
$ws.Range("B21").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-18-v.pdf)"
$ws.Range("C21").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-18-v-fasit.pdf)"

$ws.Range("A22").Value = "2018 - Høst"
$ws.Range("B22").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-18-h.pdf)"
$ws.Range("C22").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-18-h-fasit.pdf)"
$ws.Range("D22").Value = "1d, 1g"

$ws.Range("A23").Value = "2019 - Vår"
$ws.Range("B23").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-19-v.pdf)"
$ws.Range("C23").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-19-v-fasit.pdf)"
$ws.Range("D23").Value = "1c"

$ws.Range("A24").Value = "2019 - Høst"
$ws.Range("B24").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-19-h.pdf)"
$ws.Range("C24").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-19-h-fasit.pdf)"

$ws.Range("A25").Value = "2020 - Vår"
$ws.Range("B25").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-20-v.pdf)"
$ws.Range("C25").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-20-v-fasit.pdf)"

$ws.Range("A26").Value = "2020 - Høst"
$ws.Range("B26").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-20-h.pdf)"
$ws.Range("C26").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-20-h-fasit.pdf)"

# Update selection to match the saved view (last-used cell)
$ws.Range("A26").Select()

